$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 464, shifting existing rows 464:548 down to 465:549.
$ws.Rows("464:464").Insert()

# Populate the newly inserted row 464 with the new weekly record
# (same underlying values as the previous row 464 had, except for the
# date and the Origen, which reflect the new data point).
$ws.Range("A464").Value = 4
$ws.Range("B464").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C464").Value = "Los Lagos"
$ws.Range("D464").Value = 45218
$ws.Range("D464").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E464").Value = 10
$ws.Range("F464").Value = 100114014
$ws.Range("G464").Value = "Betarraga"
$ws.Range("H464").Value = "Sin especificar"
$ws.Range("I464").Value = "Primera"
$ws.Range("J464").Value = 500
$ws.Range("K464").Value = 1000
$ws.Range("L464").Value = 1000
$ws.Range("M464").Value = 1000
$ws.Range("N464").Value = "`$/paquete 5 unidades"
$ws.Range("O464").Value = "Región Metropolitana"
$ws.Range("P464").Value = 200
$ws.Range("Q464").Value = 5
$ws.Range("R464").Value = "Hortaliza"
